$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns stay text so Excel does not
# auto-convert numeric-looking strings (e.g. "1.01") into numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "90.938.21"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.123.34"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "239.55"
$ws.Range("E5").Value = "  +9.77%  "
$ws.Range("D6").Value = "630.82"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "1.05"
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("D8").Value = "0.361"
$ws.Range("E8").Value = "  -3.74%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "3.123.67"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "0.714"
$ws.Range("E11").Value = "  -7.60%  "
$ws.Range("D12").Value = "0.196"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "36.87"
$ws.Range("E13").Value = "  +5.25%  "
$ws.Range("D14").Value = "0.0000245"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").Value = "5.49"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "90.770.04"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "3.718.24"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "3.097.50"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "3.77"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "14.19"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "0.0000208"
$ws.Range("E21").Value = "  -4.14%  "
$ws.Range("D22").Value = "443.23"
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").Value = "5.56"
$ws.Range("E23").Value = "  +6.67%  "
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "5.74"
$ws.Range("E25").Value = "  -8.88%  "
$ws.Range("D26").Value = "84.12"
$ws.Range("E26").Value = "  -5.93%  "
$ws.Range("D27").Value = "12.57"
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "9.60"
$ws.Range("E30").Value = "  +3.96%  "
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("E32").Value = "  +7.84%  "
$ws.Range("D33").Value = "26.52"
$ws.Range("E33").Value = "  +6.34%  "
$ws.Range("D34").Value = "0.193"
$ws.Range("E34").Value = "  +18.61%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "510.89"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "3.76"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "0.148"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "7.19"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("D40").Value = "1.29"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").Value = "0.413"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "22.20"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").Value = "0.0839"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  +47.75%  "
$ws.Range("D46").Value = "1.91"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("D47").Value = "150.77"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").Value = "0.683"
$ws.Range("E48").Value = "  +7.41%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "45.43"
$ws.Range("E49").Value = "  +2.46%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("D51").Value = "4.48"
$ws.Range("E51").Value = "  +3.25%  "

# Restore default (unstyled) formatting so only values changed.
$dataRange.Style = "Normal"

